$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated 2D training schedule values (rows 2-6, columns B-H)
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 6
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = -3
$ws.Range("H2").Value = 34

$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = -1
$ws.Range("H3").Value = 12

$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 8
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = -5
$ws.Range("H4").Value = 56

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 7
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = -4
$ws.Range("H5").Value = 45

$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 8
$ws.Range("D6").Value = 6
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 4
$ws.Range("G6").Value = -2
$ws.Range("H6").Value = 23
